$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '35.340.26'
Set-TextValue 'D3' '1.911.23'
Set-TextValue 'E3' '  +2.84%  '
Set-TextValue 'E4' '  -0.46%  '
Set-TextValue 'D5' '245.85'
Set-TextValue 'E5' '  +2.69%  '
Set-TextValue 'E6' '  +6.18%  '
Set-TextValue 'E7' '  -0.40%  '
Set-TextValue 'D8' '41.26'
Set-TextValue 'E8' '  -2.21%  '
Set-TextValue 'D9' '0.350'
Set-TextValue 'E9' '  +6.06%  '
Set-TextValue 'D10' '52.80'
Set-TextValue 'E10' '  +12.54%  '
Set-TextValue 'D11' '0.0717'
Set-TextValue 'E11' '  +3.41%  '
Set-TextValue 'D12' '0.0994'
Set-TextValue 'E12' '  +0.52%  '
Set-TextValue 'D13' '2.188.93'
Set-TextValue 'E13' '  +2.90%  '
Set-TextValue 'D14' '12.12'
Set-TextValue 'E14' '  +5.42%  '
Set-TextValue 'E15' '  +3.78%  '
Set-TextValue 'D16' '1.905.57'
Set-TextValue 'E16' '  +2.50%  '
Set-TextValue 'D17' '4.86'
Set-TextValue 'E17' '  +3.05%  '
Set-TextValue 'D18' '35.337.94'
Set-TextValue 'D19' '72.31'
Set-TextValue 'E19' '  +3.49%  '
Set-TextValue 'D20' '0.0₃0831'
Set-TextValue 'E20' '  +4.24%  '
Set-TextValue 'D21' '239.59'
Set-TextValue 'E21' '  -0.39%  '
Set-TextValue 'D22' '12.50'
Set-TextValue 'E22' '  +2.26%  '
Set-TextValue 'D23' '4.85'
Set-TextValue 'E23' '  +2.39%  '
Set-TextValue 'E24' '  -0.40%  '
Set-TextValue 'B25' 'Toncoin'
Set-TextValue 'C25' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D25' '2.30'
Set-TextValue 'E25' '  +1.41%  '
Set-TextValue 'B26' 'PancakeSwap'
Set-TextValue 'C26' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D26' '2.35'
Set-TextValue 'E26' '  +23.04%  '
Set-TextValue 'D27' '169.67'
Set-TextValue 'E27' '  +0.65%  '
Set-TextValue 'E28' '  +6.48%  '
Set-TextValue 'D29' '18.46'
Set-TextValue 'E29' '  +4.64%  '
Set-TextValue 'E30' '  +2.15%  '
Set-TextValue 'E31' '  +3.91%  '
Set-TextValue 'E32' '  +1.89%  '
Set-TextValue 'B33' 'ImmutableX'
Set-TextValue 'C33' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D33' '0.936'
Set-TextValue 'E33' '  +14.25%  '
Set-TextValue 'B34' 'BinanceUSD'
Set-TextValue 'C34' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D34' '1.01'
Set-TextValue 'E34' '  -0.41%  '
Set-TextValue 'D35' '4.13'
Set-TextValue 'E35' '  +2.80%  '
Set-TextValue 'D36' '1.75'
Set-TextValue 'E36' '  -4.19%  '
Set-TextValue 'E37' '  +1.08%  '
Set-TextValue 'D38' '1.33'
Set-TextValue 'E38' '  +1.07%  '
Set-TextValue 'E39' '  +2.00%  '
Set-TextValue 'D40' '0.0655'
Set-TextValue 'E40' '  +8.73%  '
Set-TextValue 'E41' '  +4.10%  '
Set-TextValue 'D42' '16.36'
Set-TextValue 'E42' '  +9.34%  '
Set-TextValue 'D43' '90.15'
Set-TextValue 'E43' '  +0.22%  '
Set-TextValue 'D44' '1.338.45'
Set-TextValue 'E44' '  -0.29%  '
Set-TextValue 'E45' '  +3.33%  '
Set-TextValue 'D46' '48.10'
Set-TextValue 'E46' '  +38.22%  '
Set-TextValue 'B47' 'MXToken'
Set-TextValue 'C47' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D47' '2.79'
Set-TextValue 'E47' '  +1.96%  '
Set-TextValue 'B48' 'HuobiToken'
Set-TextValue 'C48' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D48' '2.41'
Set-TextValue 'E48' '  -0.33%  '
Set-TextValue 'D49' '6.59'
Set-TextValue 'E49' '  -0.14%  '
Set-TextValue 'D50' '2.095.69'
Set-TextValue 'E50' '  +2.75%  '
Set-TextValue 'D51' '0.0704'
Set-TextValue 'E51' '  +3.61%  '

Write-Host "Applied all crypto list updates"
